# "Tabla de servicios" - commit the data needed for the services table
# (issue #527 "tabla de servicios")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Row 4: first service row -------------------------------------------
# Column map (per header row 3): B=No. Servicio C=Fecha solicitud
# D=Fecha entrega E=Estado F=Tipo G=Encargado BE H=Encargado FE
# I=Descripcion del servicio J=Estructura JSON K=Descripcion campos

# Dates (Fecha solicitud / Fecha entrega) - 05/09/2017
$fecha = Get-Date -Year 2017 -Month 5 -Day 9 -Hour 0 -Minute 0 -Second 0
$ws.Range("C4").Value = $fecha
$ws.Range("C4").NumberFormat = "mm-dd-yy"
$ws.Range("D4").Value = $fecha
$ws.Range("D4").NumberFormat = "mm-dd-yy"

# Responsables / tipo (order matters so new shared strings line up the
# same way the authored workbook does)
$ws.Range("G4").Value = "Juan Camilo Lancheros"
$ws.Range("H4").Value = "Cristian Garzon"
$ws.Range("F4").Value = "Nuevo"

# Estructura JSON
$json = "var desbloqueoAlfaNumerico=  {`"id_jugador`":{`n                 `"id_personaje`": `"codigoAlfa`",`n  `"estado`":`"0`" }`n};"
$ws.Range("J4").Value = $json

# Descripcion del servicio
$ws.Range("I4").Value = "El servicio nos brindara el desbloueo del personaje personaje, cumpliendo si el codigo aun es valido por la fecha y si no esta desbloqueado el personaje. "

# Descripcion campos
$ws.Range("K4").Value = "id_jugador= sera el id de cada usuario en el jugo. Id_personaje= sera la informacion de el personaje que se desbloqueara. codigoAlfa = al codigo alfanumericos que se ha regalado.  Estado= Nos dira si esta desbloqueado o no este personaje."

# Wrap the long text columns so the content is readable
$ws.Range("I4").WrapText = $true
$ws.Range("J4").WrapText = $true
$ws.Range("K4").WrapText = $true

# Grow the row to fit the new wrapped content
$ws.Rows.Item(4).RowHeight = 159

# Widen the columns that now hold longer content
$ws.Columns.Item(5).ColumnWidth = 30.6640625   # E
$ws.Columns.Item(6).ColumnWidth = 28.88671875  # F
$ws.Columns.Item(7).ColumnWidth = 23.33203125  # G
$ws.Columns.Item(9).ColumnWidth = 34.33203125  # I
$ws.Columns.Item(10).ColumnWidth = 33.88671875 # J
$ws.Columns.Item(11).ColumnWidth = 45.109375   # K

# --- View housekeeping ----------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 70
$ws.Range("K8").Select()
